$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add "NA" values under the duplicate_image_filename column (column E)
# for rows 2 through 21.
$ws.Range("E2:E21").Value = "NA"
